# Natmi following Dr Hou advice
# Rewrite the Osm-Lifr sending/target cluster data to include "sCs" as an
# additional sending cluster (rows 5-7) and target cluster (already used),
# and refresh all the computed statistics columns E:T.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A:D columns (sending cluster, ligand, receptor, target cluster) per row
$rowsAD = @{
    2 = @("ECs", "Osm", "Lifr", "ECs")
    3 = @("ECs", "Osm", "Lifr", "FAPs")
    4 = @("ECs", "Osm", "Lifr", "sCs")
    5 = @("sCs", "Osm", "Lifr", "ECs")
    6 = @("sCs", "Osm", "Lifr", "FAPs")
    7 = @("sCs", "Osm", "Lifr", "sCs")
}

# E:T numeric columns per row
$rowsET = @{
    2 = @(3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729, 3, 1, 14.96835266666667, 44.905058, 0.1240053612000741, 0.1240053612000741, 327.8916792022496, 2951.025112820246, 0.123786827499404, 0.123786827499404)
    3 = @(3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729, 3, 1, 80.77474466666666, 242.324234, 0.6691786071115035, 0.6691786071115035, 1769.424281729218, 15924.81853556296, 0.6679993187645634, 0.6679993187645634)
    4 = @(3, 1, 21.90566233333334, 65.716987, 0.9982377076397728, 0.9982377076397729, 3, 1, 24.96420533333334, 74.892616, 0.2068160316884225, 0.2068160316884225, 546.8574524519992, 4921.717072067992, 0.2064515613758054, 0.2064515613758055)
    5 = @(1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098, 3, 1, 14.96835266666667, 44.905058, 0.1240053612000741, 0.1240053612000741, 0.5788611237762221, 5.209750113986, 0.0002185337006700923, 0.0002185337006700924)
    6 = @(1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098, 3, 1, 80.77474466666666, 242.324234, 0.6691786071115035, 0.6691786071115035, 3.123747850664222, 28.113730655978, 0.001179288346940013, 0.001179288346940013)
    7 = @(1, 0.3333333333333333, 0.03867233333333333, 0.116017, 0.001762292360227098, 0.001762292360227098, 3, 1, 24.96420533333334, 74.892616, 0.2068160316884225, 0.2068160316884225, 0.9654240700524445, 8.688816630472001, 0.0003644703126169923, 0.0003644703126169924)
}

foreach ($r in 2..7) {
    $ad = $rowsAD[$r]
    $ws.Range("A$r").Value = $ad[0]
    $ws.Range("B$r").Value = $ad[1]
    $ws.Range("C$r").Value = $ad[2]
    $ws.Range("D$r").Value = $ad[3]

    $et = $rowsET[$r]
    for ($i = 0; $i -lt $et.Length; $i++) {
        $col = 5 + $i  # E = column 5
        $ws.Cells.Item($r, $col).Value = $et[$i]
    }
}
